$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellOoxml {
    param(
        [object]$Cell,
        [string]$BodyXml
    )
    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
$BodyXml
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
    $Cell.Range.InsertXML($xml)
}

# Row: 3.3.2023 | 2 | "Värikoodaus duedaten mukaan"
$body1 = @"
<w:p><w:r><w:t>3.3.2023</w:t></w:r></w:p>
"@
Set-CellOoxml $t.Cell(15, 1) $body1

$body2 = @"
<w:p><w:r><w:t>2</w:t></w:r></w:p>
"@
Set-CellOoxml $t.Cell(15, 2) $body2

$body3 = @"
<w:p>
  <w:r><w:t xml:space="preserve">Värikoodaus </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>duedaten</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> mukaan</w:t></w:r>
</w:p>
"@
Set-CellOoxml $t.Cell(15, 3) $body3

# Row: 4.3.2023 | 3 | "Drag&Drop refaktoroinnin tulokseton tutkiminen"
$body4 = @"
<w:p><w:r><w:t>4.3.2023</w:t></w:r></w:p>
"@
Set-CellOoxml $t.Cell(16, 1) $body4

$body5 = @"
<w:p><w:r><w:t>3</w:t></w:r></w:p>
"@
Set-CellOoxml $t.Cell(16, 2) $body5

$body6 = @"
<w:p>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Drag&amp;Drop</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>refaktoroinnin</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> tulokseton tutkiminen</w:t></w:r>
</w:p>
"@
Set-CellOoxml $t.Cell(16, 3) $body6

# Row: 7.3.2023 | 1,5 | "Code Cleanup + Refresh nappula ja tietokannan päivitys"
$body7 = @"
<w:p><w:r><w:t>7.3.2023</w:t></w:r></w:p>
"@
Set-CellOoxml $t.Cell(17, 1) $body7

$body8 = @"
<w:p>
  <w:r><w:t>1</w:t></w:r>
  <w:r><w:t>,5</w:t></w:r>
</w:p>
"@
Set-CellOoxml $t.Cell(17, 2) $body8

$body9 = @"
<w:p>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Code</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Cleanup</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> + </w:t></w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r><w:t>Refresh</w:t></w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r><w:t xml:space="preserve"> nappula ja tietoka</w:t></w:r>
  <w:r><w:t>nnan päivitys</w:t></w:r>
</w:p>
"@
Set-CellOoxml $t.Cell(17, 3) $body9

Write-Output "done"
